$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 11000.625
$ws.Range("I62").Value = 3602.5
$ws.Range("J62").Value = 13466.667
$ws.Range("K62").Value = 3602.5
$ws.Range("L62").Value = 13466.667
$ws.Range("M62").Value = -2978.5
$ws.Range("N62").Value = -14714.667
$ws.Range("H65").Value = 11000.625
$ws.Range("I65").Value = 3602.5
$ws.Range("J65").Value = 13466.667
$ws.Range("K65").Value = 18012.5
$ws.Range("L65").Value = 67333.33499999999
$ws.Range("M65").Value = -14892.5
$ws.Range("N65").Value = -73573.33499999999
$ws.Range("H86").Value = 65227280
$ws.Range("I86").Value = 88248290
$ws.Range("J86").Value = 1095
$ws.Range("K86").Value = 88248290
$ws.Range("L86").Value = 1095
$ws.Range("M86").Value = -88247167
$ws.Range("N86").Value = -3341
$ws.Range("H89").Value = 65227280
$ws.Range("I89").Value = 88248290
$ws.Range("J89").Value = 1095
$ws.Range("K89").Value = 441241450
$ws.Range("L89").Value = 5475
$ws.Range("M89").Value = -441235834
$ws.Range("N89").Value = -16707
$ws.Range("H107").Value = 1032.0555
$ws.Range("I107").Value = 1037.6428
$ws.Range("J107").Value = 1012.5
$ws.Range("K107").Value = 1037.6428
$ws.Range("L107").Value = 1012.5
$ws.Range("M107").Value = 882.3571999999999
$ws.Range("N107").Value = -4852.5
$ws.Range("H132").Value = 3078.4055
$ws.Range("I132").Value = 2938.2942
$ws.Range("K132").Value = 8814.882599999999
$ws.Range("M132").Value = -6284.882599999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 10871200
$ws.Range("I74").Value = 951.29034
$ws.Range("K74").Value = 951.29034
$ws.Range("M74").Value = -77.29034000000001
$ws.Range("H77").Value = 10871200
$ws.Range("I77").Value = 951.29034
$ws.Range("K77").Value = 4756.4517
$ws.Range("M77").Value = -388.4516999999996
$ws.Range("H110").Value = 50748.5
$ws.Range("I110").Value = 75751.836
$ws.Range("K110").Value = 75751.836
$ws.Range("M110").Value = -73706.836
$ws.Range("H122").Value = 85219.914
$ws.Range("I122").Value = 92694.45
$ws.Range("K122").Value = 278083.35
$ws.Range("M122").Value = -275633.35
$ws.Range("H132").Value = 1151001.1
$ws.Range("J132").Value = 4053324.2
$ws.Range("L132").Value = 12159972.6
$ws.Range("N132").Value = -12165032.6

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H75").Value = 20232.46
$ws.Range("J75").Value = 27724.666
$ws.Range("L75").Value = 27724.666
$ws.Range("N75").Value = -29596.666
$ws.Range("H78").Value = 20232.46
$ws.Range("J78").Value = 27724.666
$ws.Range("L78").Value = 83173.99800000001
$ws.Range("N78").Value = -92533.99800000001
$ws.Range("H86").Value = 1888.7441
$ws.Range("I86").Value = 1865
$ws.Range("J86").Value = 2010.8572
$ws.Range("K86").Value = 1865
$ws.Range("L86").Value = 2010.8572
$ws.Range("M86").Value = -742
$ws.Range("N86").Value = -4256.8572
$ws.Range("H89").Value = 1888.7441
$ws.Range("I89").Value = 1865
$ws.Range("J89").Value = 2010.8572
$ws.Range("K89").Value = 9325
$ws.Range("L89").Value = 10054.286
$ws.Range("M89").Value = -3709
$ws.Range("N89").Value = -21286.286
$ws.Range("H134").Value = 3020.7317
$ws.Range("I134").Value = 3013.2122
$ws.Range("K134").Value = 9039.6366
$ws.Range("M134").Value = -6504.6366

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H32").Value = 0
$ws.Range("I32").Value = 0
$ws.Range("K32").Value = 0
$ws.Range("M32").ClearContents()
$ws.Range("H134").Value = 8070173.5
$ws.Range("I134").Value = 10422519
$ws.Range("J134").Value = 4989.5713
$ws.Range("K134").Value = 31267557
$ws.Range("L134").Value = 14968.7139
$ws.Range("M134").Value = -31265022
$ws.Range("N134").Value = -20038.7139

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 83333496
$ws.Range("I23").Value = 162.5
$ws.Range("J23").Value = 107143020
$ws.Range("K23").Value = 487.5
$ws.Range("L23").Value = 321429060
$ws.Range("M23").Value = -252.5
$ws.Range("N23").Value = -321429530
$ws.Range("H113").Value = 625.6667
$ws.Range("I113").Value = 602.3333
$ws.Range("J113").Value = 655.6667
$ws.Range("K113").Value = 1806.9999
$ws.Range("L113").Value = 1967.0001
$ws.Range("M113").Value = 363.0001
$ws.Range("N113").Value = -6307.0001
$ws.Range("H114").Value = 1091.9166
$ws.Range("J114").Value = 1986.6666
$ws.Range("L114").Value = 5959.9998
$ws.Range("N114").Value = -12467.9998
$ws.Range("H117").Value = 2782
$ws.Range("I117").Value = 2192.5
$ws.Range("J117").Value = 2968.158
$ws.Range("K117").Value = 6577.5
$ws.Range("L117").Value = 8904.474
$ws.Range("M117").Value = -3135.5
$ws.Range("N117").Value = -15788.474
$ws.Range("H136").Value = 3849.9443
$ws.Range("I136").Value = 1618.091
$ws.Range("J136").Value = 7357.143
$ws.Range("K136").Value = 4854.272999999999
$ws.Range("L136").Value = 22071.429
$ws.Range("M136").Value = 245.7270000000008
$ws.Range("N136").Value = -32271.429
$ws.Range("H139").Value = 305996.88
$ws.Range("I139").Value = 371978.88
$ws.Range("J139").Value = 9077.833000000001
$ws.Range("K139").Value = 1115936.64
$ws.Range("L139").Value = 27233.499
$ws.Range("M139").Value = -1110796.64
$ws.Range("N139").Value = -37513.499
$ws.Range("H141").Value = 6492.9375
$ws.Range("I141").Value = 2972.7856
$ws.Range("K141").Value = 8918.356800000001
$ws.Range("M141").Value = -3738.356800000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H112").Value = 77404.11
$ws.Range("J112").Value = 77404.11
$ws.Range("L112").Value = 77404.11
$ws.Range("N112").Value = -79620.11
$ws.Range("H122").Value = 2895
$ws.Range("I122").Value = 4082.5
$ws.Range("J122").Value = 2216.4285
$ws.Range("K122").Value = 12247.5
$ws.Range("L122").Value = 6649.2855
$ws.Range("M122").Value = -9797.5
$ws.Range("N122").Value = -11549.2855
$ws.Range("H132").Value = 33339512
$ws.Range("I132").Value = 55563300
$ws.Range("J132").Value = 3828.25
$ws.Range("K132").Value = 166689900
$ws.Range("L132").Value = 11484.75
$ws.Range("M132").Value = -166687370
$ws.Range("N132").Value = -16544.75
$ws.Range("H133").Value = 60780
$ws.Range("J133").Value = 60780
$ws.Range("L133").Value = 60780
$ws.Range("N133").Value = -70900

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 649999.5
$ws.Range("J62").Value = 300000
$ws.Range("L62").Value = 300000
$ws.Range("N62").Value = -301248
$ws.Range("H65").Value = 649999.5
$ws.Range("J65").Value = 300000
$ws.Range("L65").Value = 1500000
$ws.Range("N65").Value = -1506240
$ws.Range("H81").Value = 4646
$ws.Range("I81").Value = 7745
$ws.Range("J81").Value = 3096.5
$ws.Range("K81").Value = 15490
$ws.Range("L81").Value = 6193
$ws.Range("M81").Value = -14429
$ws.Range("N81").Value = -8315
$ws.Range("H84").Value = 4646
$ws.Range("I84").Value = 7745
$ws.Range("J84").Value = 3096.5
$ws.Range("K84").Value = 77450
$ws.Range("L84").Value = 30965
$ws.Range("M84").Value = -72146
$ws.Range("N84").Value = -41573
$ws.Range("H133").Value = 84357.5
$ws.Range("J133").Value = 84357.5
$ws.Range("L133").Value = 84357.5
$ws.Range("N133").Value = -94477.5
$ws.Range("H136").Value = 2781.4055
$ws.Range("I136").Value = 2563.9333
$ws.Range("J136").Value = 3713.4285
$ws.Range("K136").Value = 7691.7999
$ws.Range("L136").Value = 11140.2855
$ws.Range("M136").Value = -5141.7999
$ws.Range("N136").Value = -16240.2855

Write-Host "edit complete"
